$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = "上期"
$ws.Range("C3").Value = "日盘"
$ws.Range("D3").Value = "2017-02-23"
$ws.Range("E3").Value = "1.更新上期所日盘参数`n2.更新so`n麻烦检查一下主力合约`n1. 19(ok)`n2. 63(ok)"
$ws.Range("F3").Value = "passed"

$ws.Rows.Item(3).RowHeight = 119.25
